$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after the header/first data row (before old row 3).
# This pushes the former rows 3..31 down to 5..33, matching the diff's net effect
# (dimension grows from A1:T31 to A1:T33).
$ws.Rows("3:4").Insert()

# ---- New row 3 ----
$ws.Cells.Item(3,1).Value  = 1
$ws.Cells.Item(3,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(3,4).Value  = 45169
$ws.Cells.Item(3,5).Value  = 15
$ws.Cells.Item(3,6).Value  = "Fruta"
$ws.Cells.Item(3,7).Value  = 100101
$ws.Cells.Item(3,8).Value  = "Berries"
$ws.Cells.Item(3,9).Value  = 100101007
$ws.Cells.Item(3,10).Value = "Kiwi"
$ws.Cells.Item(3,11).Value = "Hayward"
$ws.Cells.Item(3,12).Value = "Primera"
$ws.Cells.Item(3,13).Value = 270
$ws.Cells.Item(3,14).Value = 27000
$ws.Cells.Item(3,15).Value = 28000
$ws.Cells.Item(3,16).Value = 27500
$ws.Cells.Item(3,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(3,18).Value = "Región de O'Higgins"
$ws.Cells.Item(3,19).Value = 1528
$ws.Cells.Item(3,20).Value = 18

# ---- New row 4 ----
$ws.Cells.Item(4,1).Value  = 1
$ws.Cells.Item(4,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(4,4).Value  = 45169
$ws.Cells.Item(4,5).Value  = 15
$ws.Cells.Item(4,6).Value  = "Fruta"
$ws.Cells.Item(4,7).Value  = 100101
$ws.Cells.Item(4,8).Value  = "Berries"
$ws.Cells.Item(4,9).Value  = 100101007
$ws.Cells.Item(4,10).Value = "Kiwi"
$ws.Cells.Item(4,11).Value = "Hayward"
$ws.Cells.Item(4,12).Value = "Segunda"
$ws.Cells.Item(4,13).Value = 150
$ws.Cells.Item(4,14).Value = 26000
$ws.Cells.Item(4,15).Value = 26000
$ws.Cells.Item(4,16).Value = 26000
$ws.Cells.Item(4,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(4,18).Value = "Región de O'Higgins"
$ws.Cells.Item(4,19).Value = 1444
$ws.Cells.Item(4,20).Value = 18
